$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.205.21"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.655.94"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").Value = "'218.15"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.5318"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "'0.2625"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "'0.06342"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "'20.40"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "'0.07827"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "'4.524"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "1.669.80"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "1.882.22"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'0.5495"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "0.0₅8158"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'65.41"
$ws.Range("D18").Value = "26.166.25"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'4.595"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'191.31"
$ws.Range("E21").Value = "  -1.51%  "
$ws.Range("D22").Value = "'10.10"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'6.022"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'145.28"
$ws.Range("E25").Value = "  +3.85%  "
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").Value = "'7.200"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'15.99"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").Value = "'1.472"
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("D30").Value = "'0.05757"
$ws.Range("E30").Value = "  -3.43%  "
$ws.Range("D31").Value = "'1.276"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").Value = "'3.555"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").Value = "'3.268"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "'1.592"
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("D35").Value = "'2.808"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'0.9514"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.424"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "'0.5759"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("D39").Value = "'0.01600"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.8521"
$ws.Range("D41").Value = "'5.799"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'1.007"
$ws.Range("E42").Value = "  +0.53%  "
$ws.Range("D43").Value = "1.044.25"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("D44").Value = "'103.93"
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("D45").Value = "1.794.75"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'56.70"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").Value = "'0.4367"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").Value = "'0.05157"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'7.843"
$ws.Range("E51").Value = "  -0.14%  "
